$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text runs (volume number and week-covering dates) ---
$ws.Range("A8").Characters(21, 2).Text = "48"
$ws.Range("C9").Characters(27, 10).Text = "11/24/2025"
$ws.Range("C9").Characters(48, 10).Text = "11/30/2025"

# --- Update crime-statistics numeric table (rows 15-30) ---
# Cells that were previously blank-placeholder text ("0") become real numbers;
# give them the same #,##0 numeric style used by sibling numeric cells.
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("C29").Value = 1
$ws.Range("C29").NumberFormat = "#,##0"
$ws.Range("F29").Value = 1
$ws.Range("F29").NumberFormat = "#,##0"
$ws.Range("C30").Value = 1
$ws.Range("C30").NumberFormat = "#,##0"
$ws.Range("F30").Value = 1
$ws.Range("F30").NumberFormat = "#,##0"

# Remaining plain numeric value updates
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 2
$ws.Range("I15").Value = 16
$ws.Range("J15").Value = 11
$ws.Range("K15").Value = 45.454545454545
$ws.Range("L15").Value = -5.882352941176
$ws.Range("M15").Value = -33.333333333333
$ws.Range("N15").Value = -78.666666666666
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -41.666666666666
$ws.Range("I16").Value = 109
$ws.Range("J16").Value = 120
$ws.Range("K16").Value = -9.166666666666
$ws.Range("L16").Value = -22.695035460992
$ws.Range("M16").Value = -65.830721003134
$ws.Range("N16").Value = -90.328305235137
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -22.222222222222
$ws.Range("F17").Value = 20
$ws.Range("H17").Value = -4.761904761904
$ws.Range("I17").Value = 223
$ws.Range("J17").Value = 278
$ws.Range("K17").Value = -19.784172661870
$ws.Range("L17").Value = -20.071684587813
$ws.Range("M17").Value = -25.418060200668
$ws.Range("N17").Value = -72.298136645962
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -36.363636363636
$ws.Range("I18").Value = 88
$ws.Range("J18").Value = 78
$ws.Range("K18").Value = 12.820512820512
$ws.Range("L18").Value = -31.250000000000
$ws.Range("M18").Value = -60.538116591928
$ws.Range("N18").Value = -88.451443569553
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -45.454545454545
$ws.Range("F19").Value = 25
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = -3.846153846153
$ws.Range("I19").Value = 258
$ws.Range("J19").Value = 258
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = -19.122257053291
$ws.Range("M19").Value = -22.522522522522
$ws.Range("N19").Value = -26.074498567335
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 200
$ws.Range("I20").Value = 76
$ws.Range("K20").Value = -12.643678160919
$ws.Range("L20").Value = -22.448979591836
$ws.Range("M20").Value = -1.298701298701
$ws.Range("N20").Value = -86.330935251798
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = -24
$ws.Range("F21").Value = 67
$ws.Range("G21").Value = 74
$ws.Range("H21").Value = -9.459459459459
$ws.Range("I21").Value = 773
$ws.Range("J21").Value = 836
$ws.Range("K21").Value = -7.535885167464
$ws.Range("L21").Value = -21.681864235055
$ws.Range("M21").Value = -40.170278637770
$ws.Range("N21").Value = -79.102460124357
$ws.Range("F22").Value = 2
$ws.Range("I22").Value = 16
$ws.Range("K22").Value = -11.111111111111
$ws.Range("L22").Value = 14.285714285714
$ws.Range("M22").Value = -38.461538461538
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 300
$ws.Range("I23").Value = 98
$ws.Range("K23").Value = 18.072289156626
$ws.Range("L23").Value = 22.500000000000
$ws.Range("M23").Value = 25.641025641025
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 5
$ws.Range("E24").Value = 240
$ws.Range("F24").Value = 63
$ws.Range("G24").Value = 57
$ws.Range("H24").Value = 10.526315789473
$ws.Range("I24").Value = 666
$ws.Range("J24").Value = 661
$ws.Range("K24").Value = 0.756429652042
$ws.Range("L24").Value = -14.285714285714
$ws.Range("M24").Value = -10.604026845637
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 11
$ws.Range("G25").Value = 7
$ws.Range("H25").Value = 57.142857142857
$ws.Range("I25").Value = 80
$ws.Range("J25").Value = 123
$ws.Range("K25").Value = -34.959349593495
$ws.Range("L25").Value = -56.043956043956
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -33.333333333333
$ws.Range("F26").Value = 35
$ws.Range("G26").Value = 18
$ws.Range("H26").Value = 94.444444444444
$ws.Range("I26").Value = 406
$ws.Range("J26").Value = 379
$ws.Range("K26").Value = 7.124010554089
$ws.Range("L26").Value = -11.739130434782
$ws.Range("M26").Value = -45.938748335552
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 22
$ws.Range("J27").Value = 17
$ws.Range("K27").Value = 29.411764705882
$ws.Range("L27").Value = -12
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 4
$ws.Range("H28").Value = 33.333333333333
$ws.Range("I28").Value = 36
$ws.Range("J28").Value = 36
$ws.Range("L28").Value = 28.571428571428
$ws.Range("H29").Value = -50
$ws.Range("I29").Value = 17
$ws.Range("K29").Value = -10.526315789473
$ws.Range("L29").Value = -5.555555555555
$ws.Range("M29").Value = -69.090909090909
$ws.Range("N29").Value = -89.759036144578
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 13
$ws.Range("K30").Value = -18.750000000000
$ws.Range("L30").Value = -7.142857142857
$ws.Range("M30").Value = -72.340425531914
$ws.Range("N30").Value = -91.333333333333
